# Auto-generated Excel COM-interop script
# Applies cell value corrections to the Adamantoise_Profits workbook sheets
# as described by the source diff (re-run of a scheduled profit-calculation pass).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 955.86664
$ws.Range("I11").Value = 955.86664
$ws.Range("K11").Value = 955.86664
$ws.Range("M11").Value = -815.86664
$ws.Range("H39").Value = 12500050
$ws.Range("I39").Value = 16666704
$ws.Range("K39").Value = 50000112
$ws.Range("M39").Value = -49999816
$ws.Range("H40").Value = 61510.33
$ws.Range("I40").Value = 602599.8
$ws.Range("K40").Value = 602599.8
$ws.Range("M40").Value = -602424.8
$ws.Range("H98").Value = 1703.2
$ws.Range("I98").Value = 1482.5
$ws.Range("K98").Value = 1482.5
$ws.Range("M98").Value = 15.5
$ws.Range("H122").Value = 1703.2
$ws.Range("I122").Value = 1482.5
$ws.Range("K122").Value = 4447.5
$ws.Range("M122").Value = -1997.5
$ws.Range("H138").Value = 3869.0244
$ws.Range("I138").Value = 1512
$ws.Range("J138").Value = 5713.6523
$ws.Range("K138").Value = 4536
$ws.Range("L138").Value = 17140.9569
$ws.Range("M138").Value = 604
$ws.Range("N138").Value = -27420.9569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2064.6667
$ws.Range("I2").Value = 1358.7693
$ws.Range("K2").Value = 1358.7693
$ws.Range("M2").Value = -1245.7693
$ws.Range("H61").Value = 3178.4375
$ws.Range("I61").Value = 2846.8572
$ws.Range("K61").Value = 2846.8572
$ws.Range("M61").Value = -2634.8572
$ws.Range("H63").Value = 120042
$ws.Range("I63").Value = 2242.3333
$ws.Range("K63").Value = 2242.3333
$ws.Range("M63").Value = -1556.3333
$ws.Range("H66").Value = 120042
$ws.Range("I66").Value = 2242.3333
$ws.Range("K66").Value = 11211.6665
$ws.Range("M66").Value = -7779.666499999999
$ws.Range("H116").Value = 2064.6667
$ws.Range("I116").Value = 1358.7693
$ws.Range("K116").Value = 1358.7693
$ws.Range("M116").Value = 935.2307000000001
$ws.Range("H132").Value = 2523.4792
$ws.Range("I132").Value = 2216.9736
$ws.Range("K132").Value = 6650.9208
$ws.Range("M132").Value = -4120.9208
$ws.Range("H136").Value = 3178.4375
$ws.Range("I136").Value = 2846.8572
$ws.Range("K136").Value = 8540.571599999999
$ws.Range("M136").Value = -5990.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2064.6667
$ws.Range("I3").Value = 1358.7693
$ws.Range("K3").Value = 1358.7693
$ws.Range("M3").Value = -1244.7693
$ws.Range("H22").Value = 328.5
$ws.Range("J22").Value = 319
$ws.Range("L22").Value = 319
$ws.Range("N22").Value = -665
$ws.Range("H86").Value = 700.3333
$ws.Range("I86").Value = 503
$ws.Range("K86").Value = 503
$ws.Range("M86").Value = 620
$ws.Range("H89").Value = 700.3333
$ws.Range("I89").Value = 503
$ws.Range("K89").Value = 2515
$ws.Range("M89").Value = 3101
$ws.Range("H92").Value = 125000
$ws.Range("J92").Value = 125000
$ws.Range("L92").Value = 125000
$ws.Range("N92").Value = -129992
$ws.Range("H117").Value = 117980
$ws.Range("J117").Value = 117980
$ws.Range("L117").Value = 117980
$ws.Range("N117").Value = -127158
$ws.Range("H134").Value = 17546572
$ws.Range("I134").Value = 2431.1
$ws.Range("K134").Value = 7293.299999999999
$ws.Range("M134").Value = -4758.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 497
$ws.Range("I22").Value = 497
$ws.Range("K22").Value = 497
$ws.Range("M22").Value = -147
$ws.Range("H31").Value = 2831.7693
$ws.Range("I31").Value = 2065.25
$ws.Range("J31").Value = 4058.2
$ws.Range("K31").Value = 2065.25
$ws.Range("L31").Value = 4058.2
$ws.Range("M31").Value = -1770.25
$ws.Range("N31").Value = -4648.2
$ws.Range("H34").Value = 2831.7693
$ws.Range("I34").Value = 2065.25
$ws.Range("J34").Value = 4058.2
$ws.Range("K34").Value = 2065.25
$ws.Range("L34").Value = 4058.2
$ws.Range("M34").Value = -1863.25
$ws.Range("N34").Value = -4462.2
$ws.Range("H58").Value = 3000.7307
$ws.Range("I58").Value = 2611
$ws.Range("J58").Value = 4058.5715
$ws.Range("K58").Value = 2611
$ws.Range("L58").Value = 4058.5715
$ws.Range("M58").Value = -2408
$ws.Range("N58").Value = -4464.5715
$ws.Range("H92").Value = 34833
$ws.Range("J92").Value = 34833
$ws.Range("L92").Value = 34833
$ws.Range("N92").Value = -39825
$ws.Range("H97").Value = 96387
$ws.Range("J97").Value = 96387
$ws.Range("L97").Value = 96387
$ws.Range("N97").Value = -98369
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2502
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 49470.094
$ws.Range("I107").Value = 60230.707
$ws.Range("K107").Value = 60230.707
$ws.Range("M107").Value = -58310.707
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9530
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 3000.7307
$ws.Range("I136").Value = 2611
$ws.Range("J136").Value = 4058.5715
$ws.Range("K136").Value = 7833
$ws.Range("L136").Value = 12175.7145
$ws.Range("M136").Value = -5283
$ws.Range("N136").Value = -17275.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 74124340
$ws.Range("I4").Value = 81324660
$ws.Range("J4").Value = 66624010
$ws.Range("K4").Value = 243973980
$ws.Range("L4").Value = 199872030
$ws.Range("M4").Value = -243973868
$ws.Range("N4").Value = -199872254
$ws.Range("H39").Value = 2945
$ws.Range("J39").Value = 4890
$ws.Range("L39").Value = 14670
$ws.Range("N39").Value = -15258
$ws.Range("H44").Value = 2726.6
$ws.Range("J44").Value = 6966.6665
$ws.Range("L44").Value = 20899.9995
$ws.Range("N44").Value = -21695.9995
$ws.Range("H69").Value = 975
$ws.Range("I69").Value = 950
$ws.Range("K69").Value = 2850
$ws.Range("M69").Value = -2039
$ws.Range("H72").Value = 975
$ws.Range("I72").Value = 950
$ws.Range("K72").Value = 8550
$ws.Range("M72").Value = -4494
$ws.Range("H132").Value = 252.8
$ws.Range("J132").Value = 200
$ws.Range("L132").Value = 1800
$ws.Range("N132").Value = -6860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 10367
$ws.Range("J92").Value = 10367
$ws.Range("L92").Value = 10367
$ws.Range("N92").Value = -14111
$ws.Range("H132").Value = 3638
$ws.Range("I132").Value = 3554.4443
$ws.Range("K132").Value = 10663.3329
$ws.Range("M132").Value = -8133.332900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7230.84
$ws.Range("I7").Value = 8463.736999999999
$ws.Range("K7").Value = 8463.736999999999
$ws.Range("M7").Value = -8351.736999999999
$ws.Range("H22").Value = 2892.543
$ws.Range("I22").Value = 2054.5
$ws.Range("J22").Value = 3779.8823
$ws.Range("K22").Value = 2054.5
$ws.Range("L22").Value = 3779.8823
$ws.Range("M22").Value = -1759.5
$ws.Range("N22").Value = -4369.8823
$ws.Range("H27").Value = 2892.543
$ws.Range("I27").Value = 2054.5
$ws.Range("J27").Value = 3779.8823
$ws.Range("K27").Value = 2054.5
$ws.Range("L27").Value = 3779.8823
$ws.Range("M27").Value = -1947.5
$ws.Range("N27").Value = -3993.8823
$ws.Range("H46").Value = 2933.9285
$ws.Range("J46").Value = 3101.923
$ws.Range("L46").Value = 3101.923
$ws.Range("N46").Value = -3477.923
$ws.Range("H55").Value = 247.94595
$ws.Range("I55").Value = 216.41176
$ws.Range("K55").Value = 216.41176
$ws.Range("M55").Value = -43.41175999999999
$ws.Range("H126").Value = 7230.84
$ws.Range("I126").Value = 8463.736999999999
$ws.Range("K126").Value = 25391.211
$ws.Range("M126").Value = -22921.211
$ws.Range("H130").Value = 112607
$ws.Range("J130").Value = 112607
$ws.Range("L130").Value = 112607
$ws.Range("N130").Value = -122647
$ws.Range("H132").Value = 3027.6155
$ws.Range("I132").Value = 2551.111
$ws.Range("J132").Value = 4099.75
$ws.Range("K132").Value = 7653.333
$ws.Range("L132").Value = 12299.25
$ws.Range("M132").Value = -5123.333
$ws.Range("N132").Value = -17359.25
$ws.Range("H136").Value = 4600.8
$ws.Range("J136").Value = 6500
$ws.Range("L136").Value = 19500
$ws.Range("N136").Value = -24600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4824.9165
$ws.Range("I62").Value = 3916.1667
$ws.Range("K62").Value = 3916.1667
$ws.Range("M62").Value = -3292.1667
$ws.Range("H65").Value = 4824.9165
$ws.Range("I65").Value = 3916.1667
$ws.Range("K65").Value = 19580.8335
$ws.Range("M65").Value = -16460.8335
$ws.Range("H131").Value = 131998
$ws.Range("J131").Value = 131998
$ws.Range("L131").Value = 131998
$ws.Range("N131").Value = -142078
$ws.Range("H132").Value = 2246.4412
$ws.Range("I132").Value = 2217.9827
$ws.Range("J132").Value = 2411.5
$ws.Range("K132").Value = 6653.9481
$ws.Range("L132").Value = 7234.5
$ws.Range("M132").Value = -4123.9481
$ws.Range("N132").Value = -12294.5
